# Applies the "inverse transformation of lines" edit:
#   - Bus: add "bus 1" and "bus 2"
#   - Load: add "demand 1" and "demand 2"
#   - Line: add "line 0-1" (bus 0 -> bus 1) and "line 1-2" (bus 1 -> bus 2)
#   - Generator: rename the diesel generator's header/bus/p_nom
#       (p_nom_max -> p_nom_min header, bus 0 -> bus 1, p_nom 200 -> 300)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Bus sheet: add two new buses
# ---------------------------------------------------------------------
$wsBus = $wb.Worksheets.Item("Bus")
$wsBus.Cells.Item(3, 1).Value = "bus 1"
$wsBus.Cells.Item(3, 2).Value = "AC"
$wsBus.Cells.Item(3, 3).Value = 2
$wsBus.Cells.Item(3, 4).Value = 0
$wsBus.Cells.Item(3, 5).Value = 0.4

$wsBus.Cells.Item(4, 1).Value = "bus 2"
$wsBus.Cells.Item(4, 2).Value = "AC"
$wsBus.Cells.Item(4, 3).Value = 4
$wsBus.Cells.Item(4, 4).Value = 0
$wsBus.Cells.Item(4, 5).Value = 0.4

$wsBus.Range("E5").Select()

# ---------------------------------------------------------------------
# Load sheet: add two new demands, one per new bus
# ---------------------------------------------------------------------
$wsLoad = $wb.Worksheets.Item("Load")
$wsLoad.Cells.Item(3, 1).Value = "demand 1"
$wsLoad.Cells.Item(3, 2).Value = "AC"
$wsLoad.Cells.Item(3, 3).Value = "bus 1"

$wsLoad.Cells.Item(4, 1).Value = "demand 2"
$wsLoad.Cells.Item(4, 2).Value = "AC"
$wsLoad.Cells.Item(4, 3).Value = "bus 2"

$wsLoad.Range("C5").Select()

# ---------------------------------------------------------------------
# Line sheet: add the two lines connecting bus0-bus1-bus2
# ---------------------------------------------------------------------
$wsLine = $wb.Worksheets.Item("Line")

# header H1 keeps the same text ("s_nom") - rewritten for consistency
$wsLine.Cells.Item(1, 8).Value = "s_nom"

$wsLine.Cells.Item(2, 1).Value = "line 0-1"
$wsLine.Cells.Item(2, 2).Value = "AC"
$wsLine.Cells.Item(2, 3).Value = "bus 0"
$wsLine.Cells.Item(2, 4).Value = "bus 1"
$wsLine.Cells.Item(2, 5).Value = "False"
$wsLine.Cells.Item(2, 6).Value = 1
$wsLine.Cells.Item(2, 7).Value = 0.1
$wsLine.Cells.Item(2, 8).Value = 100

$wsLine.Cells.Item(3, 1).Value = "line 1-2"
$wsLine.Cells.Item(3, 2).Value = "AC"
$wsLine.Cells.Item(3, 3).Value = "bus 1"
$wsLine.Cells.Item(3, 4).Value = "bus 2"
$wsLine.Cells.Item(3, 5).Value = "False"
$wsLine.Cells.Item(3, 6).Value = 1
$wsLine.Cells.Item(3, 7).Value = 0.1
$wsLine.Cells.Item(3, 8).Value = 100

$wsLine.Range("H4").Select()

# ---------------------------------------------------------------------
# Generator sheet: the diesel generator now sits on "bus 1" with a
# higher nominal power, and the p_nom_max column becomes p_nom_min
# ---------------------------------------------------------------------
$wsGen = $wb.Worksheets.Item("Generator")
$wsGen.Cells.Item(1, 7).Value = "p_nom_min"
$wsGen.Cells.Item(2, 3).Value = "bus 1"
$wsGen.Cells.Item(2, 8).Value = 300

$wsGen.Range("C3").Select()
